# Refresh the "Price" (D) and "Volume(1h)" (E) columns of the cryptos table
# with newly scraped values, per
# "Updated cryptos list on Fri Oct 27 21:38:38 UTC 2023 with GitHub Actions".
#
# Values are plain text in the sheet (not numbers), so price strings that
# would otherwise be auto-recognised as numbers (e.g. "224.56") are entered
# with a leading apostrophe to force text, matching the existing column data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "33.844.49"
$ws.Range("E2").Value = "  -0.74%  "
$ws.Range("D3").Value = "1.781.17"
$ws.Range("E3").Value = "  -0.90%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'224.56"
$ws.Range("E5").Value = "  +0.88%  "
$ws.Range("D6").Value = "'0.546"
$ws.Range("E6").Value = "  -0.97%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "'32.08"
$ws.Range("E8").Value = "  -0.35%  "
$ws.Range("D9").Value = "'0.289"
$ws.Range("E9").Value = "  +1.58%  "
$ws.Range("D10").Value = "'0.0677"
$ws.Range("E10").Value = "  -5.09%  "
$ws.Range("E11").Value = "  +1.33%  "
$ws.Range("D12").Value = "2.036.67"
$ws.Range("E12").Value = "  -0.94%  "
$ws.Range("D13").Value = "'11.27"
$ws.Range("E13").Value = "  +5.47%  "
$ws.Range("D14").Value = "1.787.09"
$ws.Range("E14").Value = "  -0.57%  "
$ws.Range("D15").Value = "33.858.53"
$ws.Range("E15").Value = "  -0.78%  "
$ws.Range("D16").Value = "'0.612"
$ws.Range("E16").Value = "  -2.92%  "
$ws.Range("D17").Value = "'4.14"
$ws.Range("E17").Value = "  -1.66%  "
$ws.Range("D18").Value = "'66.61"
$ws.Range("E18").Value = "  -2.20%  "
$ws.Range("D19").Value = "'238.79"
$ws.Range("E19").Value = "  -3.09%  "
$ws.Range("D20").Value = "0.0₃0774"
$ws.Range("E20").Value = "  -1.20%  "
$ws.Range("D21").Value = "'1.00"
$ws.Range("D22").Value = "'10.58"
$ws.Range("E22").Value = "  -2.32%  "
$ws.Range("D23").Value = "'4.01"
$ws.Range("E23").Value = "  -1.63%  "
$ws.Range("E24").Value = "  -2.58%  "
$ws.Range("D25").Value = "'160.49"
$ws.Range("E25").Value = "  +0.90%  "
$ws.Range("D26").Value = "'7.03"
$ws.Range("E26").Value = "  -0.26%  "
$ws.Range("D27").Value = "'16.09"
$ws.Range("E27").Value = "  -2.81%  "
$ws.Range("E28").Value = "  -0.15%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("E30").Value = "  +1.39%  "
$ws.Range("D31").Value = "'0.0511"
$ws.Range("E31").Value = "  -2.38%  "
$ws.Range("D32").Value = "'3.59"
$ws.Range("E32").Value = "  -3.34%  "
$ws.Range("D33").Value = "'3.51"
$ws.Range("E33").Value = "  +0.32%  "
$ws.Range("D34").Value = "'1.82"
$ws.Range("E34").Value = "  -1.47%  "
$ws.Range("D35").Value = "1.386.45"
$ws.Range("E35").Value = "  -1.80%  "
$ws.Range("D36").Value = "'0.639"
$ws.Range("E36").Value = "  -0.97%  "
$ws.Range("E37").Value = "  -2.00%  "
$ws.Range("E38").Value = "  -0.87%  "
$ws.Range("E39").Value = "  +2.35%  "
$ws.Range("D40").Value = "'2.24"
$ws.Range("E40").Value = "  +4.86%  "
$ws.Range("D41").Value = "'78.57"
$ws.Range("E41").Value = "  -2.08%  "
$ws.Range("D42").Value = "'0.911"
$ws.Range("E42").Value = "  -3.31%  "
$ws.Range("D43").Value = "'13.60"
$ws.Range("E43").Value = "  +14.00%  "
$ws.Range("E44").Value = "  -2.79%  "
$ws.Range("D45").Value = "0.0₆0142"
$ws.Range("E45").Value = "  +14.97%  "
$ws.Range("D46").Value = "'0.0505"
$ws.Range("E46").Value = "  +1.95%  "
$ws.Range("E47").Value = "  +3.13%  "
$ws.Range("D48").Value = "'107.80"
$ws.Range("E48").Value = "  +0.56%  "
$ws.Range("D49").Value = "'5.85"
$ws.Range("E49").Value = "  -1.89%  "
$ws.Range("D50").Value = "1.938.14"
$ws.Range("E50").Value = "  -1.01%  "
$ws.Range("E51").Value = "  +0.05%  "
